$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.121.64'
$ws.Range("E2").Value = '  +6.19%  '

$ws.Range("D3").Value = '3.869.83'
$ws.Range("E3").Value = '  +8.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '424.85'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.57'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.49%  '

$ws.Range("D7").Value = '3.867.42'
$ws.Range("E7").Value = '  +8.35%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.612'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.47%  '

$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.728'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.160'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.80%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000342'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +10.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.04'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.31'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.85%  '

$ws.Range("D15").Value = '4.465.65'
$ws.Range("E15").Value = '  +7.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.89'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +22.15%  '

$ws.Range("D17").Value = '3.853.55'
$ws.Range("E17").Value = '  +7.96%  '

$ws.Range("E18").Value = '  -0.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.00'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.21%  '

$ws.Range("D20").Value = '67.297.86'
$ws.Range("E20").Value = '  +6.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.09'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '415.67'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.96'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.53'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.75%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.05'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.86%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '37.74'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.97'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.26'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.31'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.17'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +35.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '731.07'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.38%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.18'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.123'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.77'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("E36").Value = '  -2.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '38.91'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.18%  '

$ws.Range("E38").Value = '  +29.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '55.45'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.52%  '

$ws.Range("D40").Value = '0.0₃0753'
$ws.Range("E40").Value = '  +21.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0465'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.41%  '

$ws.Range("E42").Value = '  +3.03%  '

$ws.Range("E43").Value = '  +0.69%  '

$ws.Range("B44").Value = 'LidoDAOToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.37'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.65%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.135'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.66%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.14'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.79%  '

$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.315'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.95%  '

$ws.Range("E48").Value = '  +2.16%  '

$ws.Range("E49").Value = '  +2.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '140.70'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.87%  '

$ws.Range("E51").Value = '  +2.16%  '
